# Updates the cryptocurrency price/volume table (columns D and E,
# rows 2-51) on the active worksheet with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new Price (D), new Volume(1h) (E) -- $null means unchanged
$changes = @(
    @(2, '27.705.73', '  -0.99%  '),
    @(3, '1.757.50', '  -0.66%  '),
    @(4, '1.002', '  -0.06%  '),
    @(5, '324.75', '  +0.91%  '),
    @(6, $null, '  -0.12%  '),
    @(7, '0.4620', '  +8.32%  '),
    @(8, '0.3609', '  -0.31%  '),
    @(9, '0.07517', '  +0.61%  '),
    @(10, $null, '  -2.75%  '),
    @(11, '1.100', '  +0.25%  '),
    @(12, $null, '  -0.04%  '),
    @(13, '20.81', '  -1.22%  '),
    @(14, '6.023', '  -1.42%  '),
    @(15, '7.126', '  -2.83%  '),
    @(16, '1.759.89', '  -1.88%  '),
    @(17, '92.38', '  -0.61%  '),
    @(18, $null, '  +0.53%  '),
    @(19, '0.06397', '  -0.20%  '),
    @(20, '1.000', '  -0.08%  '),
    @(21, '16.81', '  -2.06%  '),
    @(22, '5.814', '  -2.90%  '),
    @(23, '27.774.45', '  -0.77%  '),
    @(24, $null, '  -0.48%  '),
    @(25, '2.106', '  -0.02%  '),
    @(26, '164.54', '  +3.64%  '),
    @(27, '20.39', '  +0.35%  '),
    @(28, '1.963.91', '  -1.35%  '),
    @(29, '2.080', '  -3.92%  '),
    @(30, '126.43', '  +0.40%  '),
    @(31, '1.065', '  -7.94%  '),
    @(32, '0.09276', '  +4.06%  '),
    @(33, '3.669', '  -1.92%  '),
    @(34, '5.539', '  -1.82%  '),
    @(35, '11.91', '  -4.49%  '),
    @(36, $null, '  -0.53%  '),
    @(37, '0.2104', '  -0.37%  '),
    @(38, '0.06042', '  +0.54%  '),
    @(39, '0.6358', '  -0.13%  '),
    @(40, '4.977', '  -1.27%  '),
    @(41, '1.201', '  +1.43%  '),
    @(42, '1.378', '  -2.05%  '),
    @(43, '7.840', '  +0.02%  '),
    @(44, '13.32', '  -0.45%  '),
    @(45, '0.5918', '  -0.31%  '),
    @(46, '3.715', '  +0.44%  '),
    @(47, '123.50', '  +0.57%  '),
    @(48, '1.953', '  -2.97%  '),
    @(49, '1.150', '  -3.13%  '),
    @(50, '0.06873', $null),
    @(51, '72.42', '  -2.65%  ')
)

foreach ($row in $changes) {
    $r = $row[0]
    $dVal = $row[1]
    $eVal = $row[2]

    if ($null -ne $dVal) {
        $dCell = $ws.Cells.Item($r, 4)
        # Force text storage so numeric-looking strings (e.g. "1.002",
        # "1.000") keep their literal digits instead of being parsed as
        # numbers and normalized (which would drop trailing zeros).
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.Style = "Normal"
    }

    if ($null -ne $eVal) {
        $eCell = $ws.Cells.Item($r, 5)
        $eCell.Value = $eVal
    }
}
